$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}


$ws = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $ws.Range("A2") "-266.230831769856 - 2x_1 + 3.6141338336460276y_1 - 0.1682301438399002y_2"
Set-TextValue $ws.Range("B2") "268.730831769856"
Set-TextValue $ws.Range("D2") "0.92"
Set-TextValue $ws.Range("F2") "4.699999999999999"
Set-TextValue $ws.Range("A3") "32.09709193245777 + x_1 - 3x_2 - 0.28893058161350843y_1 + 0.9080675422138836y_2"
Set-TextValue $ws.Range("B3") "-34.09709193245777"
Set-TextValue $ws.Range("D3") "0.36"
Set-TextValue $ws.Range("E3") "8.0"
Set-TextValue $ws.Range("A4") "9.337711069418402 - 0.09193245778611647y_1 + 0.2889305816135084y_2"
Set-TextValue $ws.Range("B4") "-9.337711069418402"
Set-TextValue $ws.Range("D4") "0.49"
Set-TextValue $ws.Range("E4") "0"
Set-TextValue $ws.Range("F4") "9.3"
Set-TextValue $ws.Range("A5") "-50.92215759849905 + 0.49530956848030006y_1 - 0.8424015009380863y_2"
Set-TextValue $ws.Range("B5") "50.45215759849905"
Set-TextValue $ws.Range("D5") "0.75"
Set-TextValue $ws.Range("E5") "0"
Set-TextValue $ws.Range("F5") "0.3"

$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "51.550000000000004"
Set-TextValue $ws.Range("B2") "18.099999999999998"
Set-TextValue $ws.Range("C2") "102.2"
Set-TextValue $ws.Range("D2") "0.2"

$ws = $wb.Worksheets.Item(5)   # Vector_bf (lowercase "bf") - name lookup is case-insensitive, so use index
Set-TextValue $ws.Range("A2") "0.4525766103814899"
Set-TextValue $ws.Range("A3") "-0.6819074421513442"

$ws = $wb.Worksheets.Item(6)   # Vector_BF (uppercase "BF") - name lookup is case-insensitive, so use index
Set-TextValue $ws.Range("A2") "-6.0"
Set-TextValue $ws.Range("A3") "23.0"
Set-TextValue $ws.Range("A4") "1.8114446529080674"
Set-TextValue $ws.Range("A5") "-7.264540337711069"

$wsAlpha = $wb.Worksheets.Item("Vector_Alpha")
$wsAlpha.Range("A2").Value = 1.32
$wsAlpha.Range("A3").Value = 0.42000000000000004

Write-Host "done"
